# Generate Report for Handback
#
# The zh-cn and de-de files have now been handed back (target + handback
# files produced, in sync with en-US), so:
#   - the "Status" text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" (this is a shared string so it
#     updates the Overview sheet + both language sheets' Status column
#     automatically),
#   - the "Latest Target File" (F) and "Latest Handback File" (G) columns
#     get populated (with hyperlinks, same styling/URLs as the existing
#     Source File Name / Latest Handoff File columns for that row), and
#   - the "Latest Handback DateTime" (H) is filled in with a real
#     timestamp (previously the "0001-01-01 00:00:00" placeholder).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    Updating the cells that still hold the placeholder text updates the
#    shared string everywhere it is referenced (Overview + both sheets).
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 2) zh-cn sheet ("Latest Target File" / "Latest Handback File" columns)
#    Hyperlinks.Add() stamps the cell with the workbook's Hyperlink style
#    on its own, so no separate .Style assignment is needed.
# ---------------------------------------------------------------------
$zhcn.Range("F2").Value = "4b198396-4f25-4d29-bd36-dbe4c91b5860.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/838ad3f58dc8c674f0c67233b4b15815684b98f8/e2e/4b198396-4f25-4d29-bd36-dbe4c91b5860.md", "", "", "4b198396-4f25-4d29-bd36-dbe4c91b5860.md") | Out-Null

$zhcn.Range("G2").Value = "4b198396-4f25-4d29-bd36-dbe4c91b5860.022a54277e28870b9c08acfbc58491b14b07b471.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/782bfbc9dee4e7a14db055a0062b74fbfb9755e7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4b198396-4f25-4d29-bd36-dbe4c91b5860.022a54277e28870b9c08acfbc58491b14b07b471.zh-cn.xlf", "", "", "4b198396-4f25-4d29-bd36-dbe4c91b5860.022a54277e28870b9c08acfbc58491b14b07b471.zh-cn.xlf") | Out-Null

$zhcn.Range("F3").Value = "b07ffab4-5540-460d-9686-9f583923cf1a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/838ad3f58dc8c674f0c67233b4b15815684b98f8/e2e/b07ffab4-5540-460d-9686-9f583923cf1a.md", "", "", "b07ffab4-5540-460d-9686-9f583923cf1a.md") | Out-Null

$zhcn.Range("G3").Value = "b07ffab4-5540-460d-9686-9f583923cf1a.30512969074e4414971dd16c36f2c3c0adb067ab.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/782bfbc9dee4e7a14db055a0062b74fbfb9755e7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b07ffab4-5540-460d-9686-9f583923cf1a.30512969074e4414971dd16c36f2c3c0adb067ab.zh-cn.xlf", "", "", "b07ffab4-5540-460d-9686-9f583923cf1a.30512969074e4414971dd16c36f2c3c0adb067ab.zh-cn.xlf") | Out-Null

# zh-cn "Latest Handback DateTime" - fill in the real handback timestamp
# (was the "0001-01-01 00:00:00" placeholder)
$zhcn.Range("H2").Value = "2016-03-18 14:37:48"
$zhcn.Range("H3").Value = "2016-03-18 14:37:48"

# ---------------------------------------------------------------------
# 3) de-de sheet ("Latest Target File" / "Latest Handback File" columns)
# ---------------------------------------------------------------------
$dede.Range("F2").Value = "4b198396-4f25-4d29-bd36-dbe4c91b5860.md"
$dede.Range("F2").Style = "HyperLink"
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/838ad3f58dc8c674f0c67233b4b15815684b98f8/e2e/4b198396-4f25-4d29-bd36-dbe4c91b5860.md", "", "", "4b198396-4f25-4d29-bd36-dbe4c91b5860.md") | Out-Null

$dede.Range("G2").Value = "4b198396-4f25-4d29-bd36-dbe4c91b5860.022a54277e28870b9c08acfbc58491b14b07b471.de-de.xlf"
$dede.Range("G2").Style = "HyperLink"
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7f9c7920ade87346030e64ee8ab8b9e6b5704ad7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4b198396-4f25-4d29-bd36-dbe4c91b5860.022a54277e28870b9c08acfbc58491b14b07b471.de-de.xlf", "", "", "4b198396-4f25-4d29-bd36-dbe4c91b5860.022a54277e28870b9c08acfbc58491b14b07b471.de-de.xlf") | Out-Null

$dede.Range("F3").Value = "b07ffab4-5540-460d-9686-9f583923cf1a.md"
$dede.Range("F3").Style = "HyperLink"
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/838ad3f58dc8c674f0c67233b4b15815684b98f8/e2e/b07ffab4-5540-460d-9686-9f583923cf1a.md", "", "", "b07ffab4-5540-460d-9686-9f583923cf1a.md") | Out-Null

$dede.Range("G3").Value = "b07ffab4-5540-460d-9686-9f583923cf1a.30512969074e4414971dd16c36f2c3c0adb067ab.de-de.xlf"
$dede.Range("G3").Style = "HyperLink"
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7f9c7920ade87346030e64ee8ab8b9e6b5704ad7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b07ffab4-5540-460d-9686-9f583923cf1a.30512969074e4414971dd16c36f2c3c0adb067ab.de-de.xlf", "", "", "b07ffab4-5540-460d-9686-9f583923cf1a.30512969074e4414971dd16c36f2c3c0adb067ab.de-de.xlf") | Out-Null

# de-de "Latest Handback DateTime" - a distinct (later) handback timestamp
$dede.Range("H2").Value = "2016-03-18 14:37:54"
$dede.Range("H3").Value = "2016-03-18 14:37:54"

Write-Output "Handback report generated for zh-cn and de-de."
